$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2, 3) {
    # Column D cleared, column F gets the new value
    $ws.Cells.Item($r, 4).Value = $null
    $ws.Cells.Item($r, 6).Value = -1.405

    $ws.Cells.Item($r, 7).Value  = -0.02119535477594698   # G
    $ws.Cells.Item($r, 8).Value  = -0.02119535477594698   # H
    $ws.Cells.Item($r, 9).Value  = -0.02793018773381192   # I
    $ws.Cells.Item($r, 10).Value = -0.02793018773381192   # J
    $ws.Cells.Item($r, 11).Value = -624.8                 # K
    $ws.Cells.Item($r, 12).Value = -0.0759773818933544    # L

    $ws.Cells.Item($r, 21).Value = 486.6                  # U
    $ws.Cells.Item($r, 22).Value = 0.5203165098374679     # V
    $ws.Cells.Item($r, 23).Value = 0.8615554329840043     # W
    $ws.Cells.Item($r, 24).Value = 0.1735832634131558     # X
    $ws.Cells.Item($r, 25).Value = 0.6879721695708485     # Y
    $ws.Cells.Item($r, 26).Value = 3.979492570928282      # Z
    $ws.Cells.Item($r, 27).Value = -0.1111479745913368    # AA
    $ws.Cells.Item($r, 28).Value = 0.1119058050041281     # AB
    $ws.Cells.Item($r, 29).Value = -0.2230537795954648    # AC
    $ws.Cells.Item($r, 30).Value = 2659.5                 # AD
    $ws.Cells.Item($r, 31).Value = 3.769494145011719      # AE
    $ws.Cells.Item($r, 32).Value = 2663.269494145012      # AF
    $ws.Cells.Item($r, 33).Value = 2176.669494145012      # AG
    $ws.Cells.Item($r, 34).Value = 0.7401117331905571     # AH
    $ws.Cells.Item($r, 35).Value = 1.337741877093137      # AI
    $ws.Cells.Item($r, 36).Value = 0.6994732581942846     # AJ
    $ws.Cells.Item($r, 37).Value = 1.446994373426535      # AK
    $ws.Cells.Item($r, 38).Value = 133                    # AL
    $ws.Cells.Item($r, 39).Value = 130.06                 # AM
    $ws.Cells.Item($r, 40).Value = -47.04581638068282     # AN
    $ws.Cells.Item($r, 41).Value = -1.740601503759398     # AO
    $ws.Cells.Item($r, 42).Value = -38.50467882796766     # AP
    $ws.Cells.Item($r, 43).Value = -1.779947716438567     # AQ
}
